$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.226.51"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "2.340.88"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'545.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "2.337.81"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'23.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "2.754.03"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "60.153.42"
$ws.Range("E16").Value = "  +3.96%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "2.338.46"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D22").Value = "'313.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'63.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'7.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +8.05%  "
$ws.Range("D29").Value = "'1.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "'171.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.54%  "
$ws.Range("D32").Value = "0.0₃0729"
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("D33").Value = "'5.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("E34").Value = "  +14.09%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").Value = "'18.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'4.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("D40").Value = "'321.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.51%  "
$ws.Range("D41").Value = "'38.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").Value = "'1.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").Value = "'141.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").Value = "'3.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").Value = "'0.0945"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'19.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.26%  "
$ws.Range("D47").Value = "'0.0496"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "'0.0212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").Value = "'11.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "0.0₆0208"
$ws.Range("E51").Value = "  +16.38%  "
